$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the explanation text for the FastorSlow variable (row 7, column D)
$ws.Range("D7").Value = 'Identifies whether the question pertains to "Justin" (fast) or "Nate" (slow).'

# The longer wrapped text now needs 3 lines instead of 2, so update the
# row height to fit the new wrapped content.
$ws.Rows.Item(7).RowHeight = 51

# Update the active selection as shown in the diff
$ws.Range("D17").Select()
